$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 0.03041600386301676
$ws.Range("E3").Value = 0.02858764330546061
$ws.Range("E4").Value = 0.03171764214833577
$ws.Range("E5").Value = 0.03155186971028646
$ws.Range("E6").Value = 0.03110279242197673
$ws.Range("E7").Value = 0.03106932640075684
$ws.Range("E8").Value = 0.03144137064615885
$ws.Range("E9").Value = 0.03159577846527099
$ws.Range("E10").Value = 0.02847739855448405
$ws.Range("E11").Value = 0.02911403179168701
